# Auto-generated PowerShell Excel COM-interop script
# Applies crypto price/volume updates per commit "Updated cryptos list on Tue Jul 16 05:50:00 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.738.22"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.422.78"
$ws.Range("E3").Value = "  +2.55%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.75%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.424.23"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("E9").Value = "  +3.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.45"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").Value = "4.017.44"
$ws.Range("E13").Value = "  +2.67%  "
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.20"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "
$ws.Range("D17").Value = "63.834.68"
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "3.428.97"
$ws.Range("E18").Value = "  +3.41%  "
$ws.Range("E19").Value = "  -2.04%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.81%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "388.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.15%  "
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("B23").Value = "Polygon"
$ws.Range("C23").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.541"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.60"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.43%  "
$ws.Range("E26").Value = "  +24.25%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.85%  "
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("E30").Value = "  +8.04%  "
$ws.Range("E31").Value = "  +1.94%  "
$ws.Range("E32").Value = "  +4.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.14%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  +4.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.47"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("E38").Value = "  -1.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0780"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.06%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.84"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.30%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.02"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.901.02"
$ws.Range("E42").Value = "  +2.00%  "
$ws.Range("E43").Value = "  +1.84%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.767"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.56"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +7.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.07"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.70%  "
$ws.Range("E49").Value = "  +21.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.51"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.78%  "
$ws.Range("E51").Value = "  +5.92%  "
